$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Normalized Data")
$ws.Range("A8").Value = "test"
